$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) of the last header cell onto the new header cell
# before setting the values, so the new column matches the existing header
# styling (bold font, thin border, centered alignment) exactly.
$ws.Range("G1").Copy($ws.Range("H1"))

$ws.Cells.Item(1, 8).Value = "consequents_length"
$ws.Cells.Item(2, 8).Value = 1
